$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.205607476635514
$ws.Range("C2").Value = 0.4626168224299065
$ws.Range("J2").Value = 0.009345794392523364
$ws.Range("P2").Value = 0.2009345794392523
$ws.Range("S2").Value = 0.1214953271028037
$ws.Range("C3").Value = 0.0297029702970297
$ws.Range("J3").Value = 0.0396039603960396
$ws.Range("P3").Value = 0.6831683168316832
$ws.Range("S3").Value = 0.2475247524752475
$ws.Range("J4").Value = 0.08571428571428572
$ws.Range("P4").Value = 0.6
$ws.Range("S4").Value = 0.3142857142857143
$ws.Range("B6").Value = 0.06451612903225806
$ws.Range("D6").Value = 0.01075268817204301
$ws.Range("F6").Value = 0.04301075268817205
$ws.Range("J6").Value = 0.2419354838709677
$ws.Range("O6").Value = 0.01075268817204301
$ws.Range("Q6").Value = 0.1182795698924731
$ws.Range("R6").Value = 0.05913978494623656
$ws.Range("S6").Value = 0.4516129032258064
$ws.Range("B7").Value = 0.07792207792207792
$ws.Range("D7").Value = 0.01948051948051948
$ws.Range("F7").Value = 0.07142857142857142
$ws.Range("J7").Value = 0.1428571428571428
$ws.Range("O7").Value = 0.03246753246753246
$ws.Range("Q7").Value = 0.1883116883116883
$ws.Range("R7").Value = 0.06493506493506493
$ws.Range("S7").Value = 0.4025974025974026
$ws.Range("B8").Value = 0.053475935828877
$ws.Range("D8").Value = 0.0106951871657754
$ws.Range("F8").Value = 0.06149732620320856
$ws.Range("J8").Value = 0.09625668449197861
$ws.Range("O8").Value = 0.02406417112299465
$ws.Range("Q8").Value = 0.2112299465240642
$ws.Range("R8").Value = 0.09090909090909091
$ws.Range("S8").Value = 0.4518716577540107
$ws.Range("B9").Value = 0.08205128205128205
$ws.Range("D9").Value = 0.01538461538461539
$ws.Range("F9").Value = 0.05128205128205128
$ws.Range("J9").Value = 0.1025641025641026
$ws.Range("O9").Value = 0.03589743589743589
$ws.Range("Q9").Value = 0.1948717948717949
$ws.Range("R9").Value = 0.05641025641025641
$ws.Range("S9").Value = 0.4615384615384616
$ws.Range("B10").Value = 0.1004140786749482
$ws.Range("D10").Value = 0.02380952380952381
$ws.Range("F10").Value = 0.08488612836438923
$ws.Range("J10").Value = 0.119047619047619
$ws.Range("O10").Value = 0.008281573498964804
$ws.Range("Q10").Value = 0.1966873706004141
$ws.Range("R10").Value = 0.07453416149068323
$ws.Range("S10").Value = 0.3923395445134575
$ws.Range("G11").Value = 0.1398305084745763
$ws.Range("J11").Value = 0.06779661016949153
$ws.Range("K11").Value = 0.1652542372881356
$ws.Range("L11").Value = 0.597457627118644
$ws.Range("S11").Value = 0.02966101694915254
$ws.Range("G12").Value = 0.7551020408163265
$ws.Range("J12").Value = 0.1360544217687075
$ws.Range("K12").Value = 0.006802721088435374
$ws.Range("L12").Value = 0.04081632653061224
$ws.Range("S12").Value = 0.06122448979591837
$ws.Range("F13").Value = 0.03333333333333333
$ws.Range("G13").Value = 0.5666666666666667
$ws.Range("J13").Value = 0.3666666666666666
$ws.Range("S13").Value = 0.03333333333333333
$ws.Range("F15").Value = 0.01036269430051814
$ws.Range("H15").Value = 0.1450777202072539
$ws.Range("I15").Value = 0.07253886010362694
$ws.Range("J15").Value = 0.38860103626943
$ws.Range("K15").Value = 0.05699481865284974
$ws.Range("M15").Value = 0.0155440414507772
$ws.Range("O15").Value = 0.05699481865284974
$ws.Range("S15").Value = 0.2538860103626943
$ws.Range("F16").Value = 0.05511811023622047
$ws.Range("H16").Value = 0.1889763779527559
$ws.Range("I16").Value = 0.1023622047244094
$ws.Range("J16").Value = 0.3228346456692913
$ws.Range("K16").Value = 0.1181102362204724
$ws.Range("M16").Value = 0.01574803149606299
$ws.Range("O16").Value = 0.04724409448818898
$ws.Range("S16").Value = 0.1496062992125984
$ws.Range("F17").Value = 0.0113314447592068
$ws.Range("H17").Value = 0.2152974504249292
$ws.Range("I17").Value = 0.1019830028328612
$ws.Range("J17").Value = 0.3597733711048159
$ws.Range("K17").Value = 0.1076487252124646
$ws.Range("M17").Value = 0.0113314447592068
$ws.Range("N17").Value = 0.0028328611898017
$ws.Range("O17").Value = 0.05099150141643059
$ws.Range("S17").Value = 0.1388101983002833
$ws.Range("F18").Value = 0.007352941176470588
$ws.Range("H18").Value = 0.1764705882352941
$ws.Range("I18").Value = 0.1029411764705882
$ws.Range("J18").Value = 0.4411764705882353
$ws.Range("K18").Value = 0.06617647058823529
$ws.Range("M18").Value = 0.01470588235294118
$ws.Range("O18").Value = 0.09558823529411764
$ws.Range("S18").Value = 0.09558823529411764
$ws.Range("F19").Value = 0.01531531531531532
$ws.Range("H19").Value = 0.2063063063063063
$ws.Range("I19").Value = 0.1045045045045045
$ws.Range("J19").Value = 0.345945945945946
$ws.Range("K19").Value = 0.1054054054054054
$ws.Range("M19").Value = 0.01711711711711712
$ws.Range("O19").Value = 0.08468468468468468
$ws.Range("S19").Value = 0.1207207207207207
